$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9987039566040039
$ws.Range("B1").Value = 1.834736466407776
$ws.Range("C1").Value = 1.976743698120117
$ws.Range("D1").Value = 2.052458047866821
$ws.Range("E1").Value = 1.394486784934998
